$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.013.01'
$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("D3").Value = '2.749.41'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '351.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.22'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.545'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.578'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.95'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.11%  '
$ws.Range("E11").Value = '  +3.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0830'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.14%  '
$ws.Range("D15").Value = '3.169.86'
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").Value = '2.733.13'
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.917'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = '50.921.95'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("D22").Value = '0.0₃0952'
$ws.Range("E22").Value = '  -2.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '262.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.70'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.161'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '51.76'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0437'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0824'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.18'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("E41").Value = '  -1.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.68%  '
$ws.Range("E44").Value = '  -2.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").Value = '2.079.36'
$ws.Range("E46").Value = '  +1.22%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.914'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.65%  '
$ws.Range("E51").Value = '  +5.33%  '
